# Update model_list.xlsx for "no MPA" and "outside MPA" scenarios.
# Adds ten new model rows (42-51) on Sheet1, describing Model41, Model42,
# Model44-Model50 runs that have MPA Size "None" (no-MPA / outside-MPA
# control runs), and updates the active selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------
# Columns: A=Model Name, B=Climate, C=MPA Size, D=Evolution?, E=Movement Rate
# Populate column A first (in row order) so the new "ModelNN" shared
# strings are minted before "None", matching authoring order.
$ws.Range("A42").Value = "Model41"
$ws.Range("A43").Value = "Model42"
$ws.Range("A44").Value = "Model44"
$ws.Range("A45").Value = "Model44"
$ws.Range("A46").Value = "Model45"
$ws.Range("A47").Value = "Model46"
$ws.Range("A48").Value = "Model47"
$ws.Range("A49").Value = "Model48"
$ws.Range("A50").Value = "Model49"
$ws.Range("A51").Value = "Model50"

$ws.Range("B42").Value = "Null"
$ws.Range("B43").Value = "Null"
$ws.Range("B44").Value = "Mean"
$ws.Range("B45").Value = "Mean"
$ws.Range("B46").Value = "ENSO"
$ws.Range("B47").Value = "ENSO"
$ws.Range("B48").Value = "Shock"
$ws.Range("B49").Value = "Shock"
$ws.Range("B50").Value = "Mean Shock"
$ws.Range("B51").Value = "Mean Shock"

$ws.Range("D42").Value = "Yes"
$ws.Range("D43").Value = "No"
$ws.Range("D44").Value = "Yes"
$ws.Range("D45").Value = "No"
$ws.Range("D46").Value = "Yes"
$ws.Range("D47").Value = "No"
$ws.Range("D48").Value = "Yes"
$ws.Range("D49").Value = "No"
$ws.Range("D50").Value = "Yes"
$ws.Range("D51").Value = "No"

$ws.Range("E42").Value = "1-3 grids"
$ws.Range("E43").Value = "1-3 grids"
$ws.Range("E44").Value = "1-3 grids"
$ws.Range("E45").Value = "1-3 grids"
$ws.Range("E46").Value = "1-3 grids"
$ws.Range("E47").Value = "1-3 grids"
$ws.Range("E48").Value = "1-3 grids"
$ws.Range("E49").Value = "1-3 grids"
$ws.Range("E50").Value = "1-3 grids"
$ws.Range("E51").Value = "1-3 grids"

# MPA Size column ("None") is written last of the new text so the new
# shared string lands after all ten "ModelNN" strings.
$ws.Range("C42").Value = "None"
$ws.Range("C43").Value = "None"
$ws.Range("C44").Value = "None"
$ws.Range("C45").Value = "None"
$ws.Range("C46").Value = "None"
$ws.Range("C47").Value = "None"
$ws.Range("C48").Value = "None"
$ws.Range("C49").Value = "None"
$ws.Range("C50").Value = "None"
$ws.Range("C51").Value = "None"

# --- View state ------------------------------------------------------------
$ws.Range("D40:D51").Select() | Out-Null
